$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Create the new "metadata" sheet by duplicating "data" (keeps identical sheet
# boilerplate: sheetPr/sheetFormatPr/pageMargins/sheetView) then wipe its contents.
$dataSheet.Copy($null, $dataSheet)
$metaSheet = $wb.Worksheets.Item(2)
$metaSheet.Name = "metadata"
$metaSheet.Cells.Clear()

# --- Header row ---
$metaSheet.Cells.Item(1, 2).Value = "data_name"
$metaSheet.Cells.Item(1, 3).Value = "data_id"
$metaSheet.Cells.Item(1, 4).Value = "data_version"
$metaSheet.Cells.Item(1, 5).Value = "data_version_created"
$metaSheet.Cells.Item(1, 6).Value = "panel_query_time"
$metaSheet.Cells.Item(1, 7).Value = "panel_get_request"

# --- Data row ---
$metaSheet.Cells.Item(2, 1).Value = 0
$metaSheet.Cells.Item(2, 2).Value = "Unexplained paediatric onset end-stage renal disease"
$metaSheet.Cells.Item(2, 3).Value = 678

# "1.20" must stay textual (not be coerced to the number 1.2): force Text format
# before assigning, then drop back to the Normal style so no stray formatting
# remains on the cell.
$metaSheet.Cells.Item(2, 4).NumberFormat = "@"
$metaSheet.Cells.Item(2, 4).Value = "1.20"
$metaSheet.Cells.Item(2, 4).Style = "Normal"

$metaSheet.Cells.Item(2, 5).Value = "2021-08-03T14:41:51.253372Z"
$metaSheet.Cells.Item(2, 6).Value = "2021-10-05 14:23:04.615398"
$metaSheet.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/678/?format=json"

# Reuse the existing header/index cell styles (bold+border for the header row,
# the index-column style for A2) from "data" instead of synthesizing new ones.
$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)

# Restore "data" as the active sheet (the new sheet was active right after Copy).
$dataSheet.Activate()

# --- Refresh the "time_taken" (column F) timestamps on the "data" sheet to the
#     re-run query time. ---
$newTimestamps = @(
    "2021-10-05 14:23:04.618560",
    "2021-10-05 14:23:04.618568",
    "2021-10-05 14:23:04.618572",
    "2021-10-05 14:23:04.618574",
    "2021-10-05 14:23:04.618577",
    "2021-10-05 14:23:04.618580",
    "2021-10-05 14:23:04.618582",
    "2021-10-05 14:23:04.618585",
    "2021-10-05 14:23:04.618588",
    "2021-10-05 14:23:04.618590",
    "2021-10-05 14:23:04.618593",
    "2021-10-05 14:23:04.618595",
    "2021-10-05 14:23:04.618598",
    "2021-10-05 14:23:04.618600",
    "2021-10-05 14:23:04.618603",
    "2021-10-05 14:23:04.618606",
    "2021-10-05 14:23:04.618609",
    "2021-10-05 14:23:04.618611",
    "2021-10-05 14:23:04.618614",
    "2021-10-05 14:23:04.618617",
    "2021-10-05 14:23:04.618619",
    "2021-10-05 14:23:04.618622",
    "2021-10-05 14:23:04.618624",
    "2021-10-05 14:23:04.618627",
    "2021-10-05 14:23:04.618631",
    "2021-10-05 14:23:04.618634",
    "2021-10-05 14:23:04.618636",
    "2021-10-05 14:23:04.618639",
    "2021-10-05 14:23:04.618641",
    "2021-10-05 14:23:04.618644",
    "2021-10-05 14:23:04.618669",
    "2021-10-05 14:23:04.618673",
    "2021-10-05 14:23:04.618677",
    "2021-10-05 14:23:04.618680",
    "2021-10-05 14:23:04.618683",
    "2021-10-05 14:23:04.618685",
    "2021-10-05 14:23:04.618688",
    "2021-10-05 14:23:04.618690",
    "2021-10-05 14:23:04.618693",
    "2021-10-05 14:23:04.618696",
    "2021-10-05 14:23:04.618699",
    "2021-10-05 14:23:04.618702",
    "2021-10-05 14:23:04.618704",
    "2021-10-05 14:23:04.618707",
    "2021-10-05 14:23:04.618710",
    "2021-10-05 14:23:04.618712",
    "2021-10-05 14:23:04.618715",
    "2021-10-05 14:23:04.618718",
    "2021-10-05 14:23:04.618720",
    "2021-10-05 14:23:04.618723",
    "2021-10-05 14:23:04.618726",
    "2021-10-05 14:23:04.618728",
    "2021-10-05 14:23:04.618731",
    "2021-10-05 14:23:04.618734",
    "2021-10-05 14:23:04.618737",
    "2021-10-05 14:23:04.618739",
    "2021-10-05 14:23:04.618742",
    "2021-10-05 14:23:04.618745",
    "2021-10-05 14:23:04.618747",
    "2021-10-05 14:23:04.618750",
    "2021-10-05 14:23:04.618753",
    "2021-10-05 14:23:04.618755",
    "2021-10-05 14:23:04.618758",
    "2021-10-05 14:23:04.618761",
    "2021-10-05 14:23:04.618766",
    "2021-10-05 14:23:04.618769",
    "2021-10-05 14:23:04.618772",
    "2021-10-05 14:23:04.618774",
    "2021-10-05 14:23:04.618777",
    "2021-10-05 14:23:04.618780",
    "2021-10-05 14:23:04.618782",
    "2021-10-05 14:23:04.618785",
    "2021-10-05 14:23:04.618787",
    "2021-10-05 14:23:04.618790",
    "2021-10-05 14:23:04.618793",
    "2021-10-05 14:23:04.618796",
    "2021-10-05 14:23:04.618800",
    "2021-10-05 14:23:04.618804",
    "2021-10-05 14:23:04.618806",
    "2021-10-05 14:23:04.618809",
    "2021-10-05 14:23:04.618812",
    "2021-10-05 14:23:04.618814",
    "2021-10-05 14:23:04.618817",
    "2021-10-05 14:23:04.618819",
    "2021-10-05 14:23:04.618822",
    "2021-10-05 14:23:04.618824",
    "2021-10-05 14:23:04.618827",
    "2021-10-05 14:23:04.618829",
    "2021-10-05 14:23:04.618832",
    "2021-10-05 14:23:04.618835",
    "2021-10-05 14:23:04.618837",
    "2021-10-05 14:23:04.618840",
    "2021-10-05 14:23:04.618844",
    "2021-10-05 14:23:04.618847",
    "2021-10-05 14:23:04.618850",
    "2021-10-05 14:23:04.618852",
    "2021-10-05 14:23:04.618855",
    "2021-10-05 14:23:04.618857",
    "2021-10-05 14:23:04.618860",
    "2021-10-05 14:23:04.618863",
    "2021-10-05 14:23:04.618866",
    "2021-10-05 14:23:04.618868",
    "2021-10-05 14:23:04.618871",
    "2021-10-05 14:23:04.618873",
    "2021-10-05 14:23:04.618876",
    "2021-10-05 14:23:04.618878",
    "2021-10-05 14:23:04.618881",
    "2021-10-05 14:23:04.618884",
    "2021-10-05 14:23:04.618888",
    "2021-10-05 14:23:04.618891",
    "2021-10-05 14:23:04.618894",
    "2021-10-05 14:23:04.618897",
    "2021-10-05 14:23:04.618900",
    "2021-10-05 14:23:04.618902",
    "2021-10-05 14:23:04.618905",
    "2021-10-05 14:23:04.618907",
    "2021-10-05 14:23:04.618910",
    "2021-10-05 14:23:04.618913",
    "2021-10-05 14:23:04.618915",
    "2021-10-05 14:23:04.618918",
    "2021-10-05 14:23:04.618920",
    "2021-10-05 14:23:04.618923",
    "2021-10-05 14:23:04.618925",
    "2021-10-05 14:23:04.618928",
    "2021-10-05 14:23:04.618931",
    "2021-10-05 14:23:04.618933",
    "2021-10-05 14:23:04.618936",
    "2021-10-05 14:23:04.618938",
    "2021-10-05 14:23:04.618943",
    "2021-10-05 14:23:04.618946",
    "2021-10-05 14:23:04.618949",
    "2021-10-05 14:23:04.618952",
    "2021-10-05 14:23:04.618955",
    "2021-10-05 14:23:04.618957",
    "2021-10-05 14:23:04.618960",
    "2021-10-05 14:23:04.618963",
    "2021-10-05 14:23:04.618965",
    "2021-10-05 14:23:04.618968",
    "2021-10-05 14:23:04.618970",
    "2021-10-05 14:23:04.618973",
    "2021-10-05 14:23:04.618975",
    "2021-10-05 14:23:04.618978",
    "2021-10-05 14:23:04.618981",
    "2021-10-05 14:23:04.618983",
    "2021-10-05 14:23:04.618986",
    "2021-10-05 14:23:04.618988",
    "2021-10-05 14:23:04.618991",
    "2021-10-05 14:23:04.618994",
    "2021-10-05 14:23:04.618997",
    "2021-10-05 14:23:04.619000",
    "2021-10-05 14:23:04.619002",
    "2021-10-05 14:23:04.619005",
    "2021-10-05 14:23:04.619007",
    "2021-10-05 14:23:04.619010",
    "2021-10-05 14:23:04.619013",
    "2021-10-05 14:23:04.619015",
    "2021-10-05 14:23:04.619018",
    "2021-10-05 14:23:04.619020",
    "2021-10-05 14:23:04.619023",
    "2021-10-05 14:23:04.619026",
    "2021-10-05 14:23:04.619028",
    "2021-10-05 14:23:04.619031",
    "2021-10-05 14:23:04.619033",
    "2021-10-05 14:23:04.619036",
    "2021-10-05 14:23:04.619038",
    "2021-10-05 14:23:04.619041",
    "2021-10-05 14:23:04.619044",
    "2021-10-05 14:23:04.619046",
    "2021-10-05 14:23:04.619049",
    "2021-10-05 14:23:04.619051",
    "2021-10-05 14:23:04.619054",
    "2021-10-05 14:23:04.619056",
    "2021-10-05 14:23:04.619060",
    "2021-10-05 14:23:04.619063",
    "2021-10-05 14:23:04.619066",
    "2021-10-05 14:23:04.619069",
    "2021-10-05 14:23:04.619071",
    "2021-10-05 14:23:04.619074",
    "2021-10-05 14:23:04.619077",
    "2021-10-05 14:23:04.619080",
    "2021-10-05 14:23:04.619082",
    "2021-10-05 14:23:04.619085",
    "2021-10-05 14:23:04.619087",
    "2021-10-05 14:23:04.619090",
    "2021-10-05 14:23:04.619093",
    "2021-10-05 14:23:04.619095",
    "2021-10-05 14:23:04.619098",
    "2021-10-05 14:23:04.619101",
    "2021-10-05 14:23:04.619103",
    "2021-10-05 14:23:04.619106",
    "2021-10-05 14:23:04.619108",
    "2021-10-05 14:23:04.619111",
    "2021-10-05 14:23:04.619114",
    "2021-10-05 14:23:04.619116",
    "2021-10-05 14:23:04.619119",
    "2021-10-05 14:23:04.619122",
    "2021-10-05 14:23:04.619124",
    "2021-10-05 14:23:04.619127",
    "2021-10-05 14:23:04.619130",
    "2021-10-05 14:23:04.619132",
    "2021-10-05 14:23:04.619135",
    "2021-10-05 14:23:04.619138",
    "2021-10-05 14:23:04.619141",
    "2021-10-05 14:23:04.619143",
    "2021-10-05 14:23:04.619146",
    "2021-10-05 14:23:04.619149",
    "2021-10-05 14:23:04.619151",
    "2021-10-05 14:23:04.619154",
    "2021-10-05 14:23:04.619156",
    "2021-10-05 14:23:04.619159",
    "2021-10-05 14:23:04.619162",
    "2021-10-05 14:23:04.619164",
    "2021-10-05 14:23:04.619167",
    "2021-10-05 14:23:04.619170",
    "2021-10-05 14:23:04.619172",
    "2021-10-05 14:23:04.619175",
    "2021-10-05 14:23:04.619177",
    "2021-10-05 14:23:04.619180",
    "2021-10-05 14:23:04.619182",
    "2021-10-05 14:23:04.619185",
    "2021-10-05 14:23:04.619188",
    "2021-10-05 14:23:04.619190",
    "2021-10-05 14:23:04.619193",
    "2021-10-05 14:23:04.619195",
    "2021-10-05 14:23:04.619198",
    "2021-10-05 14:23:04.619201",
    "2021-10-05 14:23:04.619203",
    "2021-10-05 14:23:04.619206",
    "2021-10-05 14:23:04.619209",
    "2021-10-05 14:23:04.619211",
    "2021-10-05 14:23:04.619214"
)

for ($i = 0; $i -lt $newTimestamps.Length; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimestamps[$i]
}
